$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 2).Value = "احمد"
$ws.Cells.Item($row, 3).Value = "'1"
$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 3"
$ws.Cells.Item($row, 6).Value = "C3"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٠١:٣٠ م"
